$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Item_Description (column H) text tweaks
$ws.Range("H8").Value2  = "green(raw) petroleum coke (in bulk)"
$ws.Range("H10").Value2 = "green(raw) petroleum coke (in bulk)"
$ws.Range("H11").Value2 = "green(raw) petroleum coke (in bulk)"
$ws.Range("H40").Value2 = "calcined petroleum coke"
$ws.Range("H44").Value2 = "petroleum coke (graphitized)"
$ws.Range("H46").Value2 = "calcined petroleum coke"

# Importer_Name (column P) tweak
$ws.Range("P34").Value2 = "brakes india private limited"

# Row 12 recalculated USD figures
$ws.Range("T12").Value2 = 2.2023
$ws.Range("U12").Value2 = 2862.9648
$ws.Range("V12").Value2 = 3.4846

# Row 15
$ws.Range("T15").Value2 = 6.3818
$ws.Range("U15").Value2 = 6381553.8694

# Row 16
$ws.Range("T16").Value2 = 6.3818
$ws.Range("U16").Value2 = 12763107.4665

# Row 17
$ws.Range("T17").Value2 = 6.3818
$ws.Range("U17").Value2 = 15953883.9929

# Row 18
$ws.Range("T18").Value2 = 6.3818
$ws.Range("U18").Value2 = 1749471.0113

# Row 19
$ws.Range("T19").Value2 = 6.3818
$ws.Range("U19").Value2 = 3190776.7986

# Row 31
$ws.Range("T31").Value2 = 2.1837
$ws.Range("U31").Value2 = 3493.902
$ws.Range("V31").Value2 = 3.4846

# Row 34
$ws.Range("T34").Value2 = 84.9866
$ws.Range("U34").Value2 = 849841.9051
$ws.Range("V34").Value2 = 0.8075

# Row 35
$ws.Range("T35").Value2 = 83.0891
$ws.Range("U35").Value2 = 332356.4061
$ws.Range("V35").Value2 = 0.8075

# Row 36
$ws.Range("T36").Value2 = 83.0891
$ws.Range("U36").Value2 = 332356.4061
$ws.Range("V36").Value2 = 0.8075

# Row 37
$ws.Range("T37").Value2 = 671.8729
$ws.Range("U37").Value2 = 134373.5613
$ws.Range("V37").Value2 = 8.228199999999999

# Row 38
$ws.Range("T38").Value2 = 79928.6032
$ws.Range("U38").Value2 = 1918286.4776
$ws.Range("V38").Value2 = 984.9489

# Row 43
$ws.Range("T43").Value2 = 262351.3154
$ws.Range("U43").Value2 = 1311756.5536
$ws.Range("V43").Value2 = 2657.6136
